$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The retrieval function lost two stations ("Gerrit van der Veenstraat" and
# the VUmc main entrance). Fix the VUmc station name to its full name.
$ws.Range("A2").Value = "VUmc Hoofdingang"
